# Expense List Section - add the new expense rows under the existing
# "Rent" row, and push the original first expense row (Rent / 300 /
# 2025-03-12) down to become the final row of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Category, Amount, Date-serial) that will occupy rows 2-7.
# Row 7 reproduces the data that used to live in row 2 before the edit.
$category = @("Rent", "Food", "Travel", "Transport", "Rent", "Rent")
$amount   = @(5000,   3000,   1000,     2000,        3000,   300)
$date     = @(46025.229537037034, 46015.229537037034, 46015.229537037034, 46015.229537037034, 46014.229537037034, 45728.229537037034)

# The existing C2 cell already carries the date number format (style index
# referenced by s="1" in the sheet XML) - copy it down to the new date
# cells so no duplicate style gets minted.
$ws.Range("C2").Copy() | Out-Null

for ($i = 0; $i -lt $category.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $category[$i]
    $ws.Cells.Item($row, 2).Value = $amount[$i]
    $ws.Cells.Item($row, 3).Value = $date[$i]
    if ($row -gt 2) {
        $ws.Cells.Item($row, 3).PasteSpecial(-4122) | Out-Null
    }
}
